# Course-material update:
#  1. Remove the two "Example" filler slides that were dropped from the
#     "reading/cleaning" deck (the data.gov CSV-import walkthrough and the
#     data.gov xls-import walkthrough).
#  2. Refresh the cached "last edited" date field (datetimeFigureOut) that
#     is shown on the slide master / all slide layouts / the notes master
#     from 12/09/2019 to 17/09/2019.

$p = $ppt.ActivePresentation

# --- 1. Delete the two "Example" slides -----------------------------------
# Slide #10 = "Example" / import from data.gov (csv) walkthrough.
$p.Slides.Item(10).Delete()
# After that deletion the old slide #13 ("Example" / xls from data.gov)
# has shifted down to position #12.
$p.Slides.Item(12).Delete()

# --- 2. Update the cached date placeholder text ---------------------------
function Update-DatePlaceholder {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "12/09/2019") {
                $shp.TextFrame.TextRange.Text = "17/09/2019"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes
